# Update cryptos list (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 42 and 43 swap order (dogwifhat <-> Filecoin) plus updated price/volume
$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").Value = "'4.42"
$ws.Range("E42").Value = "  -6.17%  "

$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").Value = "'2.53"
$ws.Range("E43").Value = "  -12.42%  "

# Price / Volume(1h) refresh for remaining rows
$ws.Range("D2").Value = "65.627.09"
$ws.Range("E2").Value = "  -3.71%  "
$ws.Range("D3").Value = "3.474.55"
$ws.Range("E3").Value = "  -0.60%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'581.70"
$ws.Range("E5").Value = "  -2.34%  "
$ws.Range("D6").Value = "'171.93"
$ws.Range("E6").Value = "  -6.02%  "
$ws.Range("E8").Value = "  -3.88%  "
$ws.Range("D9").Value = "3.473.97"
$ws.Range("E9").Value = "  -0.50%  "
$ws.Range("E10").Value = "  -7.96%  "
$ws.Range("E11").Value = "  -2.64%  "
$ws.Range("E12").Value = "  -5.00%  "
$ws.Range("D13").Value = "4.081.84"
$ws.Range("E13").Value = "  -0.28%  "
$ws.Range("E14").Value = "  +0.22%  "
$ws.Range("D15").Value = "'29.85"
$ws.Range("E15").Value = "  -7.57%  "
$ws.Range("D16").Value = "65.731.73"
$ws.Range("E16").Value = "  -3.36%  "
$ws.Range("E17").Value = "  -4.14%  "
$ws.Range("D18").Value = "3.482.50"
$ws.Range("E18").Value = "  -0.22%  "
$ws.Range("D19").Value = "'5.93"
$ws.Range("E19").Value = "  -4.57%  "
$ws.Range("D20").Value = "'13.91"
$ws.Range("E20").Value = "  -2.17%  "
$ws.Range("D21").Value = "'366.55"
$ws.Range("E21").Value = "  -7.56%  "
$ws.Range("D22").Value = "'7.77"
$ws.Range("E22").Value = "  -2.81%  "
$ws.Range("D23").Value = "'1.00"
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("D24").Value = "'0.539"
$ws.Range("E24").Value = "  -0.39%  "
$ws.Range("D25").Value = "'72.16"
$ws.Range("E25").Value = "  -0.19%  "
$ws.Range("D26").Value = "'0.0000122"
$ws.Range("E26").Value = "  -1.29%  "
$ws.Range("D27").Value = "'9.71"
$ws.Range("E27").Value = "  -7.46%  "
$ws.Range("E28").Value = "  +0.64%  "
$ws.Range("E29").Value = "  +0.22%  "
$ws.Range("D30").Value = "'24.11"
$ws.Range("E30").Value = "  +1.98%  "
$ws.Range("D31").Value = "'5.77"
$ws.Range("E31").Value = "  -6.23%  "
$ws.Range("E32").Value = "  -3.80%  "
$ws.Range("D33").Value = "'0.999"
$ws.Range("E33").Value = "  -0.07%  "
$ws.Range("D34").Value = "'1.28"
$ws.Range("E34").Value = "  -9.11%  "
$ws.Range("D35").Value = "'7.05"
$ws.Range("E35").Value = "  -4.42%  "
$ws.Range("E36").Value = "  -1.96%  "
$ws.Range("D37").Value = "'29.53"
$ws.Range("E37").Value = "  +11.75%  "
$ws.Range("D38").Value = "'159.61"
$ws.Range("E38").Value = "  -1.40%  "
$ws.Range("E39").Value = "  -1.38%  "
$ws.Range("D40").Value = "'1.78"
$ws.Range("E40").Value = "  -5.42%  "
$ws.Range("D41").Value = "2.803.23"
$ws.Range("E41").Value = "  +1.67%  "
$ws.Range("D44").Value = "'6.30"
$ws.Range("E44").Value = "  -6.38%  "
$ws.Range("E45").Value = "  -4.97%  "
$ws.Range("E46").Value = "  -4.40%  "
$ws.Range("D47").Value = "'24.23"
$ws.Range("E47").Value = "  -8.44%  "
$ws.Range("E48").Value = "  -4.21%  "
$ws.Range("D49").Value = "'305.17"
$ws.Range("E49").Value = "  -8.02%  "
$ws.Range("D50").Value = "'0.823"
$ws.Range("E50").Value = "  -3.30%  "
$ws.Range("E51").Value = "  -4.35%  "
